$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used for the observation core fields that move together with a record.
$cols = @("A","B","E","F","G","H","Q","R")

# Snapshot current ("before") values for the rows involved in the reshuffle so that
# writes to one row never clobber data we still need to read from another row.
$rows = @(5,6,7,8,9,10,11,13)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowData
}

# Also snapshot the special per-row cells that differ between the "Tretåig hackspett"
# observation (row 9, before) and the "Granticka" observation (row 11, before):
# row 9 (before) has L9/M9/AC9 populated and lacks J9/AF9,
# row 11 (before) has J11/AF11 populated (empty placeholders) and lacks L11/M11/AC11.
$m9 = $ws.Range("M9").Value()

# Mapping of destination row -> source row (which "before" record ends up there).
$mapping = @{
    5  = 7
    6  = 5
    7  = 13
    8  = 10
    9  = 11
    10 = 8
    11 = 9
    13 = 6
}

# Write the core fields for every destination row from the snapshotted source row.
foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $data = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $data[$c]
    }
}

# Row 9 now holds the record that used to be in row 11 ("Granticka" / Porodaedalea
# chrysoloma). That record uses empty placeholder cells J and AF, and has no L/M/AC.
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("J9").Style = "Normal"
$ws.Range("AF9").Style = "Normal"

# Row 11 now holds the record that used to be in row 9 ("Tretåig hackspett" /
# Picoides tridactylus). That record has L/M/AC populated and lacks J/AF.
$ws.Range("J11").ClearContents()
$ws.Range("AF11").ClearContents()
$ws.Range("L11").Style = "Normal"
$ws.Range("M11").Value = $m9
$ws.Range("AC11").Value = "Skalad gran"
